$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50 already holds the same date ("2025/10/02") stored as plain text.
# Copy it down so A51 stays text instead of Excel auto-converting the
# "yyyy/mm/dd" string into a date serial number + date-formatted style.
$ws.Range("A50").Copy($ws.Range("A51"))

$ws.Range("B51").Value = "木"
$ws.Range("C51").Value = 16
$ws.Range("D51").Value = 3
